$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1421.625
$ws.Range("I2").Value = 874.8
$ws.Range("J2").Value = 2333
$ws.Range("K2").Value = 874.8
$ws.Range("L2").Value = 2333
$ws.Range("M2").Value = -761.8
$ws.Range("N2").Value = -2559

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 135.6
$ws.Range("I6").Value = 135.6
$ws.Range("K6").Value = 406.8
$ws.Range("M6").Value = -294.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 2740.625
$ws.Range("I8").Value = 2740.625
$ws.Range("K8").Value = 8221.875
$ws.Range("M8").Value = -8082.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 858511.2
$ws.Range("I9").Value = 219.6923
$ws.Range("K9").Value = 219.6923
$ws.Range("M9").Value = -50.69229999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3851
$ws.Range("I40").Value = 3601.4
$ws.Range("J40").Value = 4085
$ws.Range("K40").Value = 3601.4
$ws.Range("L40").Value = 4085
$ws.Range("M40").Value = -3426.4
$ws.Range("N40").Value = -4435

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 6120.1665
$ws.Range("I70").Value = 4181.5
$ws.Range("J70").Value = 8058.8335
$ws.Range("K70").Value = 12544.5
$ws.Range("L70").Value = 24176.5005
$ws.Range("M70").Value = -12274.5
$ws.Range("N70").Value = -24716.5005

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 6120.1665
$ws.Range("I73").Value = 4181.5
$ws.Range("J73").Value = 8058.8335
$ws.Range("K73").Value = 12544.5
$ws.Range("L73").Value = 24176.5005
$ws.Range("M73").Value = -11608.5
$ws.Range("N73").Value = -26048.5005

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1768.7273
$ws.Range("I80").Value = 1500.3334
$ws.Range("J80").Value = 1869.375
$ws.Range("K80").Value = 4501.0002
$ws.Range("L80").Value = 5608.125
$ws.Range("M80").Value = -3503.0002
$ws.Range("N80").Value = -7604.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1768.7273
$ws.Range("I83").Value = 1500.3334
$ws.Range("J83").Value = 1869.375
$ws.Range("K83").Value = 13503.0006
$ws.Range("L83").Value = 16824.375
$ws.Range("M83").Value = -8511.000599999999
$ws.Range("N83").Value = -26808.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2785.6875
$ws.Range("I86").Value = 2285.625
$ws.Range("J86").Value = 3285.75
$ws.Range("K86").Value = 2285.625
$ws.Range("L86").Value = 3285.75
$ws.Range("M86").Value = -1162.625
$ws.Range("N86").Value = -5531.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 2785.6875
$ws.Range("I89").Value = 2285.625
$ws.Range("J89").Value = 3285.75
$ws.Range("K89").Value = 11428.125
$ws.Range("L89").Value = 16428.75
$ws.Range("M89").Value = -5812.125
$ws.Range("N89").Value = -27660.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 4229.636
$ws.Range("I100").Value = 1807.8
$ws.Range("K100").Value = 1807.8
$ws.Range("M100").Value = -1266.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4304.4653
$ws.Range("I32").Value = 4397.1577
$ws.Range("K32").Value = 4397.1577
$ws.Range("M32").Value = -4110.1577

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 43222
$ws.Range("J109").Value = 43222
$ws.Range("L109").Value = 43222
$ws.Range("N109").Value = -45996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3017.7368
$ws.Range("I132").Value = 2517.889
$ws.Range("K132").Value = 7553.667
$ws.Range("M132").Value = -5023.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 95376
$ws.Range("J108").Value = 95376
$ws.Range("L108").Value = 95376
$ws.Range("N108").Value = -103056

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 106372.29
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 32271.773
$ws.Range("J9").Value = 32271.773
$ws.Range("L9").Value = 32271.773
$ws.Range("N9").Value = -32607.773

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 4280.3335
$ws.Range("I17").Value = 4665.5
$ws.Range("K17").Value = 4665.5
$ws.Range("M17").Value = -4491.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 11524.083
$ws.Range("I62").Value = 6111.857
$ws.Range("J62").Value = 19101.2
$ws.Range("K62").Value = 6111.857
$ws.Range("L62").Value = 19101.2
$ws.Range("M62").Value = -5487.857
$ws.Range("N62").Value = -20349.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 11524.083
$ws.Range("I65").Value = 6111.857
$ws.Range("J65").Value = 19101.2
$ws.Range("K65").Value = 30559.285
$ws.Range("L65").Value = 95506
$ws.Range("M65").Value = -27439.285
$ws.Range("N65").Value = -101746

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 131226.6
$ws.Range("I68").Value = 157895.28
$ws.Range("J68").Value = 68999.664
$ws.Range("K68").Value = 157895.28
$ws.Range("L68").Value = 68999.664
$ws.Range("M68").Value = -157146.28
$ws.Range("N68").Value = -70497.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 131226.6
$ws.Range("I71").Value = 157895.28
$ws.Range("J71").Value = 68999.664
$ws.Range("K71").Value = 473685.84
$ws.Range("L71").Value = 206998.992
$ws.Range("M71").Value = -469941.84
$ws.Range("N71").Value = -214486.992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 771.2308
$ws.Range("I107").Value = 696.125
$ws.Range("K107").Value = 696.125
$ws.Range("M107").Value = 1223.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3073.6743
$ws.Range("I134").Value = 2195.3125
$ws.Range("K134").Value = 6585.9375
$ws.Range("M134").Value = -4050.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2478.7666
$ws.Range("J107").Value = 3375.65
$ws.Range("L107").Value = 10126.95
$ws.Range("N107").Value = -13966.95

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 34860
$ws.Range("I122").Value = 34860
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 313740
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -311290
$ws.Range("N122").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2006.1333
$ws.Range("I137").Value = 1731
$ws.Range("J137").Value = 2246.875
$ws.Range("K137").Value = 5193
$ws.Range("L137").Value = 6740.625
$ws.Range("M137").Value = -93
$ws.Range("N137").Value = -16940.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 13109.333
$ws.Range("J19").Value = 13498
$ws.Range("L19").Value = 13498
$ws.Range("N19").Value = -14074

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 15000
$ws.Range("I21").Value = 10000
$ws.Range("J21").Value = 20000
$ws.Range("K21").Value = 10000
$ws.Range("L21").Value = 20000
$ws.Range("M21").Value = -9827
$ws.Range("N21").Value = -20346

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value = 15000
$ws.Range("I30").Value = 10000
$ws.Range("J30").Value = 20000
$ws.Range("K30").Value = 10000
$ws.Range("L30").Value = 20000
$ws.Range("M30").Value = -9895
$ws.Range("N30").Value = -20210

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 314999.34
$ws.Range("J34").Value = 314999.34
$ws.Range("L34").Value = 314999.34
$ws.Range("N34").Value = -315535.34

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16249.75
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 16249.75
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 16249.75
$ws.Range("M70").Value = $null
$ws.Range("N70").Value = -16789.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 16249.75
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 16249.75
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 16249.75
$ws.Range("M73").Value = $null
$ws.Range("N73").Value = -18121.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H76").Value = 314999.34
$ws.Range("J76").Value = 314999.34
$ws.Range("L76").Value = 314999.34
$ws.Range("N76").Value = -315629.34

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H79").Value = 314999.34
$ws.Range("J79").Value = 314999.34
$ws.Range("L79").Value = 314999.34
$ws.Range("N79").Value = -317183.34

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 19996
$ws.Range("J110").Value = 19996
$ws.Range("L110").Value = 19996
$ws.Range("N110").Value = -28176

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3285.4546
$ws.Range("I126").Value = 2412
$ws.Range("K126").Value = 7236
$ws.Range("M126").Value = -4766

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 81603.16
$ws.Range("I7").Value = 81603.16
$ws.Range("K7").Value = 81603.16
$ws.Range("M7").Value = -81491.16

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1838.6
$ws.Range("I61").Value = 1723.25
$ws.Range("J61").Value = 2300
$ws.Range("K61").Value = 1723.25
$ws.Range("L61").Value = 2300
$ws.Range("M61").Value = -1521.25
$ws.Range("N61").Value = -2704

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3833.1667
$ws.Range("I68").Value = 4000
$ws.Range("J68").Value = 3499.5
$ws.Range("K68").Value = 4000
$ws.Range("L68").Value = 3499.5
$ws.Range("M68").Value = -3251
$ws.Range("N68").Value = -4997.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3833.1667
$ws.Range("I71").Value = 4000
$ws.Range("J71").Value = 3499.5
$ws.Range("K71").Value = 20000
$ws.Range("L71").Value = 17497.5
$ws.Range("M71").Value = -16256
$ws.Range("N71").Value = -24985.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 27779944
$ws.Range("I100").Value = 83334830
$ws.Range("J100").Value = 2498.3333
$ws.Range("K100").Value = 83334830
$ws.Range("L100").Value = 2498.3333
$ws.Range("M100").Value = -83334289
$ws.Range("N100").Value = -3580.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1838.6
$ws.Range("I113").Value = 1723.25
$ws.Range("J113").Value = 2300
$ws.Range("K113").Value = 1723.25
$ws.Range("L113").Value = 2300
$ws.Range("M113").Value = 446.75
$ws.Range("N113").Value = -6640

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 81603.16
$ws.Range("I126").Value = 81603.16
$ws.Range("K126").Value = 244809.48
$ws.Range("M126").Value = -242339.48

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3407.147
$ws.Range("I136").Value = 2463.5518
$ws.Range("K136").Value = 7390.655400000001
$ws.Range("M136").Value = -4840.655400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 37682.6
$ws.Range("J20").Value = 42603.25
$ws.Range("L20").Value = 42603.25
$ws.Range("N20").Value = -43083.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 97086
$ws.Range("J116").Value = 97086
$ws.Range("L116").Value = 97086
$ws.Range("N116").Value = -106264

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1989.9048
$ws.Range("I136").Value = 1146.7368
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 3440.2104
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -890.2103999999999
$ws.Range("N136").Value = -35100
